$d = $word.ActiveDocument

# 1) Title paragraph: split "SPRINT 2 - Daily Scrums" so "Scrums" is in its
#    own spell-checked run (w:proofErr spellStart/spellEnd wrap it).
$p1 = $d.Paragraphs(1)
$p1.Range.InsertXML(@'
<w:p><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">SPRINT 2 - Daily </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>Scrums</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@)

# 2) Second paragraph: split "Daily Scrum - Día 1" so "Día" is in its own
#    spell-checked run.
$p2 = $d.Paragraphs(2)
$p2.Range.InsertXML(@'
<w:p><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">Daily Scrum - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>Día</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve"> 1</w:t></w:r></w:p>
'@)

# 3) Append the new "Daily Scrum - Día 2" section at the end of the body,
#    right before the closing section properties.
$end = $d.Content
$end.Collapse(0)
$end.InsertXML(@'
<w:p/><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Daily</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Scrum - Día 2</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Facilitador: Jorge Samuel Solano Dorantes (Scrum </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Master</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t>Alán Osmar Peña Polo (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Product</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Owner</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> / </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Developer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Frontend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Backend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t>¿Qué hice ayer?</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Implementé exitosamente la validación de diferencia visual entre </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sprites</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (T1.1)</w:t></w:r></w:p><w:p><w:r><w:t>El algoritmo compara bit a bit y calcula el porcentaje de diferencia</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Creé el componente de configuración de metas con radio </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>buttons</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> funcionales (T1.2)</w:t></w:r></w:p><w:p><w:r><w:t>Los inputs numéricos tienen validación de rangos (1-99 obstáculos, 1-255 segundos)</w:t></w:r></w:p><w:p><w:r><w:t>Implementé el sistema de mensajes de error descriptivos (T1.4)</w:t></w:r></w:p><w:p><w:r><w:t>Ahora los usuarios ven mensajes específicos según el tipo de error</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>¿Qué haré hoy?</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Actualizar el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>backend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Flask</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> para recibir configuración completa: personaje + obstáculo + metas (T2.1)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Modificar el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>endpoint</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> /api/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>send_config</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> para validar la nueva estructura JSON</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Agregar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>feedback</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> visual detallado en el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>frontend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (T2.3)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Implementar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>spinner</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de carga, mensajes de progreso y animaciones</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>Realizar pruebas de la nueva estructura JSON con mi compañero</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Impedimentos detectados:</w:t></w:r></w:p><w:p><w:r><w:t>IMPEDIMENTO: Necesito coordinar con mi compañero el formato exacto del JSON extendido que incluye las metas. ¿El PIC puede procesar campos adicionales sin problemas de memoria?</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Jorge Samuel Solano Dorantes (Scrum </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Master</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> / </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Developer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - Embebido)</w:t></w:r></w:p><w:p><w:r><w:t>¿Qué hice ayer?</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Implementé la estructura de datos para metas en el PIC (T1.3)</w:t></w:r></w:p><w:p><w:r><w:t>Optimicé el uso de memoria, las nuevas variables ocupan solo 3 bytes adicionales</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Verifiqué que aún tengo suficiente RAM disponible para el resto de </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>funcionalidades</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p/><w:p><w:r><w:t>¿Qué haré hoy?</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Implementar validación de JSON en el PIC según las acciones de mejora del Sprint 1 (T2.2)</w:t></w:r></w:p><w:p><w:r><w:t>Validar estructura JSON: verificar que existan los campos requeridos</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Validar rangos de valores: 0x00-0xFF para bytes de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sprites</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, 0-1 para </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>meta_type</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, 1-255 para </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>meta_value</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>Implementar respuestas de error específicas con códigos numéricos</w:t></w:r></w:p><w:p><w:r><w:t>Procesar y almacenar configuración de meta recibida (T2.4)</w:t></w:r></w:p><w:p><w:r><w:t>Coordinar con mi compañero para definir el esquema JSON final con metas</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Impedimentos detectados:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">IMPEDIMENTO: Necesito definir con mi compañero el formato JSON exacto y los códigos de error numéricos para mantener consistencia entre </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>backend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> y embebido.</w:t></w:r></w:p><w:p/><w:p/>
'@)
